# The commit "Automatic update of files" swaps the data of several pairs
# of observation rows in the sheet (the underlying records kept their row
# position but the record contents were exchanged between two rows).
#
# Pairs that swap (by row number): (2,3) (16,18) (17,19) (28,29) (30,31)
# (32,33) (61,62). For most pairs only columns A,B,E,F,G,H,Q,R differ; for
# the (61,62) pair the age/stage/activity/method columns (K,L,M,N) and the
# public-comment column (AC) also move between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cells {
    param(
        [int]$Row1,
        [int]$Row2,
        [string[]]$Columns
    )

    foreach ($col in $Columns) {
        $addr1 = "$col$Row1"
        $addr2 = "$col$Row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

$standardColumns = @("A", "B", "E", "F", "G", "H", "Q", "R")
$standardPairs = @(
    @(2, 3),
    @(16, 18),
    @(17, 19),
    @(28, 29),
    @(30, 31),
    @(32, 33)
)

foreach ($pair in $standardPairs) {
    Swap-Cells $pair[0] $pair[1] $standardColumns
}

# Row 61 / 62 swap extra columns K, L, M, N, AC on top of the standard set.
$extendedColumns = @("A", "B", "E", "F", "G", "H", "K", "L", "M", "N", "Q", "R", "AC")
Swap-Cells 61 62 $extendedColumns
